$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Style work: the two empty header cells that sit to the right of the
#    "original" value in each merged B1:D1 / E1:G1 band ("C1"/"D1" and
#    "F1"/"G1") need a lighter border (no header font/centering) instead
#    of the bold, fully-boxed header style: a top+bottom rule for the
#    inner one, plus a right rule for the one that closes the group.
#
#    Doing `.Borders(...).LineStyle = ...` repeatedly, cell-by-cell,
#    against this host makes the later cells fork off a *new* style
#    instead of reusing the one already produced for an earlier cell
#    with an identical border (a host quirk) - so instead we build the
#    two target look-and-feels once on a scratch sheet and fan them out
#    with Copy / PasteSpecial(xlPasteFormats), which reuses one xf per
#    distinct look. The scratch sheet is removed again afterwards.
# ---------------------------------------------------------------------
$scratchSheet = $wb.Worksheets.Add()

# Template A -> border: top thin, bottom thin, no left, no right
$tA = $scratchSheet.Range("A1")
$tA.Borders(8).LineStyle = 1       # xlEdgeTop    = thin
$tA.Borders(9).LineStyle = 1       # xlEdgeBottom = thin
$tA.Borders(7).LineStyle = -4142   # xlEdgeLeft   = none
$tA.Borders(10).LineStyle = -4142  # xlEdgeRight  = none

# Template B -> border: top thin, bottom thin, right thin, no left.
# Set top+bottom before left/right: this host keeps every distinct
# border combination it ever sees (even fleeting intermediate ones) in
# the saved styles.xml, so reaching the top+bottom-only combination
# first (already used by template A) and only then adding the right
# edge avoids registering a throwaway "right-only" border along the way.
$tB = $scratchSheet.Range("B1")
$tB.Borders(8).LineStyle = 1       # xlEdgeTop    = thin
$tB.Borders(9).LineStyle = 1       # xlEdgeBottom = thin
$tB.Borders(7).LineStyle = -4142   # xlEdgeLeft   = none
$tB.Borders(10).LineStyle = 1      # xlEdgeRight  = thin

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

$tA.Copy()
$ws1.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$tB.Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.DisplayAlerts = $false
$scratchSheet.Delete() | Out-Null

# Deleting the scratch sheet shifts sheet indices under the hood, so the
# worksheet handles grabbed earlier are no longer reliable - re-resolve
# them by name before doing any further editing.
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---------------------------------------------------------------------
# 2) Anonymize the "fedcore" column header -> "approach"
# ---------------------------------------------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---------------------------------------------------------------------
# 3) Normalize negative-zero "change" values to plain zero
# ---------------------------------------------------------------------
$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("D12").Value = 0

# ---------------------------------------------------------------------
# 4) Drop the stray empty inline-string cell left over in row 5
# ---------------------------------------------------------------------
$ws2.Range("G5").ClearContents()
